$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate existing labels/values from Spanish to English and update values
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Karl"
$ws.Range("A2").Value = "Age"
$ws.Range("A3").Value = "Gender"
$ws.Range("B3").Value = "Man"
$ws.Range("A4").Value = "Height"
$ws.Range("B4").Value = "1.78 m"
$ws.Range("A5").Value = "Weight"
$ws.Range("B5").Value = "72 kg"

# Add new row 6 with uppercase name info, styled like the other rows (bordered)
# but with a new blue font color
$ws.Range("A6").Value = "Uppercase name"
$ws.Range("B6").Value = "KARL"

# Reuse the existing bordered style (copy formats from row 1) then recolor the font to blue (#0070C0)
$ws.Range("A1:B1").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)
$ws.Range("A6:B6").Font.Color = 12611584

$ws.Range("B4").Select()
